# Reversi board: "Refactored possibilities marker. Added markings for all
# directions."
#
# The sheet marks legal-move "possibility" squares with a shared string
# "B" (Black) or "W" (White) letter; conditional formatting (already in
# the workbook) then colours the cell based on that letter. Around the
# white disc at F6, markers previously only covered three of the
# directions (E5/F5/G5 diagonal+vertical neighbours). This change marks
# every direction: it adds the missing up/down markers (F4 and G7), drops
# the now-redundant G5 marker, and flips the marker sitting on F6 itself
# from "W" to "B".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "B" possibility markers for the previously-unmarked directions.
$ws.Range("F4").Value = "B"
$ws.Range("G7").Value = "B"

# F6's own marker flips from "W" to "B" (keeps its highlighted style).
$ws.Range("F6").Value = "B"

# G5 is no longer marked as a possibility.
$ws.Range("G5").ClearContents()

# Cursor ends up on F6 after the edits.
$ws.Range("F6").Select() | Out-Null
